$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 857.7273
$ws.Range("I8").Value = 54.375
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 163.125
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -24.125
$ws.Range("N8").Value = -9278
$ws.Range("H9").Value = 9839
$ws.Range("J9").Value = 1133.1666
$ws.Range("L9").Value = 1133.1666
$ws.Range("N9").Value = -1471.1666
$ws.Range("H52").Value = 2125
$ws.Range("I52").Value = 250
$ws.Range("K52").Value = 750
$ws.Range("M52").Value = -590
$ws.Range("H80").Value = 460.3846
$ws.Range("J80").Value = 521
$ws.Range("L80").Value = 1563
$ws.Range("N80").Value = -3559
$ws.Range("H83").Value = 460.3846
$ws.Range("J83").Value = 521
$ws.Range("L83").Value = 4689
$ws.Range("N83").Value = -14673
$ws.Range("H98").Value = 3557.3333
$ws.Range("I98").Value = 3742.2
$ws.Range("K98").Value = 3742.2
$ws.Range("M98").Value = -2244.2
$ws.Range("H103").Value = 568.4286
$ws.Range("I103").Value = 591.6667
$ws.Range("K103").Value = 1775.0001
$ws.Range("M103").Value = -1189.0001
$ws.Range("H111").Value = 12195.143
$ws.Range("I111").Value = 6998.7144
$ws.Range("K111").Value = 20996.1432
$ws.Range("M111").Value = -17929.1432
$ws.Range("H122").Value = 3557.3333
$ws.Range("I122").Value = 3742.2
$ws.Range("K122").Value = 11226.6
$ws.Range("M122").Value = -8776.599999999999
$ws.Range("H125").Value = 3957.25
$ws.Range("I125").Value = 5250
$ws.Range("J125").Value = 2664.5
$ws.Range("K125").Value = 47250
$ws.Range("L125").Value = 23980.5
$ws.Range("M125").Value = -44790
$ws.Range("N125").Value = -28900.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1297.45
$ws.Range("I2").Value = 1213.2941
$ws.Range("K2").Value = 1213.2941
$ws.Range("M2").Value = -1100.2941
$ws.Range("H16").Value = 3486.0908
$ws.Range("I16").Value = 2378.4285
$ws.Range("J16").Value = 5424.5
$ws.Range("K16").Value = 2378.4285
$ws.Range("L16").Value = 5424.5
$ws.Range("M16").Value = -2091.4285
$ws.Range("N16").Value = -5998.5
$ws.Range("H19").Value = 12006.714
$ws.Range("I19").Value = 13984.5
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 13984.5
$ws.Range("L19").Value = 140
$ws.Range("M19").Value = -13755.5
$ws.Range("N19").Value = -598
$ws.Range("H74").Value = 12021.723
$ws.Range("I74").Value = 2106.8572
$ws.Range("K74").Value = 2106.8572
$ws.Range("M74").Value = -1232.8572
$ws.Range("H77").Value = 12021.723
$ws.Range("I77").Value = 2106.8572
$ws.Range("K77").Value = 10534.286
$ws.Range("M77").Value = -6166.286
$ws.Range("H102").Value = 2270.5667
$ws.Range("I102").Value = 1874.9565
$ws.Range("J102").Value = 3570.4285
$ws.Range("K102").Value = 1874.9565
$ws.Range("L102").Value = 3570.4285
$ws.Range("M102").Value = -252.9565
$ws.Range("N102").Value = -6814.4285
$ws.Range("H116").Value = 1297.45
$ws.Range("I116").Value = 1213.2941
$ws.Range("K116").Value = 1213.2941
$ws.Range("M116").Value = 1080.7059
$ws.Range("H122").Value = 1577.8948
$ws.Range("I122").Value = 1278.5385
$ws.Range("K122").Value = 3835.6155
$ws.Range("M122").Value = -1385.6155
$ws.Range("H124").Value = 39999.668
$ws.Range("J124").Value = 39999.668
$ws.Range("L124").Value = 39999.668
$ws.Range("N124").Value = -49819.668
$ws.Range("H132").Value = 3089.1177
$ws.Range("I132").Value = 2794.926
$ws.Range("K132").Value = 8384.778
$ws.Range("M132").Value = -5854.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1297.45
$ws.Range("I3").Value = 1213.2941
$ws.Range("K3").Value = 1213.2941
$ws.Range("M3").Value = -1099.2941
$ws.Range("H107").Value = 951.3333
$ws.Range("I107").Value = 871.4400000000001
$ws.Range("K107").Value = 871.4400000000001
$ws.Range("M107").Value = 1048.56

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 316.18182
$ws.Range("I22").Value = 275.44446
$ws.Range("K22").Value = 275.44446
$ws.Range("M22").Value = 74.55554000000001
$ws.Range("H31").Value = 33229.938
$ws.Range("I31").Value = 38651.445
$ws.Range("K31").Value = 38651.445
$ws.Range("M31").Value = -38356.445
$ws.Range("H34").Value = 33229.938
$ws.Range("I34").Value = 38651.445
$ws.Range("K34").Value = 38651.445
$ws.Range("M34").Value = -38449.445
$ws.Range("H62").Value = 5878.2144
$ws.Range("J62").Value = 5050.75
$ws.Range("L62").Value = 5050.75
$ws.Range("N62").Value = -6298.75
$ws.Range("H65").Value = 5878.2144
$ws.Range("J65").Value = 5050.75
$ws.Range("L65").Value = 25253.75
$ws.Range("N65").Value = -31493.75
$ws.Range("H122").Value = 1910.5714
$ws.Range("I122").Value = 1571.2858
$ws.Range("J122").Value = 2249.8572
$ws.Range("K122").Value = 4713.857400000001
$ws.Range("L122").Value = 6749.571599999999
$ws.Range("M122").Value = -2263.857400000001
$ws.Range("N122").Value = -11649.5716
$ws.Range("H132").Value = 3207.0571
$ws.Range("I132").Value = 3012.7778
$ws.Range("J132").Value = 3862.75
$ws.Range("K132").Value = 9038.3334
$ws.Range("L132").Value = 11588.25
$ws.Range("M132").Value = -6508.3334
$ws.Range("N132").Value = -16648.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 345.8846
$ws.Range("I23").Value = 284.63635
$ws.Range("J23").Value = 390.8
$ws.Range("K23").Value = 853.90905
$ws.Range("L23").Value = 1172.4
$ws.Range("M23").Value = -618.90905
$ws.Range("N23").Value = -1642.4
$ws.Range("H81").Value = 4822
$ws.Range("I81").Value = 1513
$ws.Range("J81").Value = 5649.25
$ws.Range("K81").Value = 4539
$ws.Range("L81").Value = 16947.75
$ws.Range("M81").Value = -3416
$ws.Range("N81").Value = -19193.75
$ws.Range("H84").Value = 4822
$ws.Range("I84").Value = 1513
$ws.Range("J84").Value = 5649.25
$ws.Range("K84").Value = 13617
$ws.Range("L84").Value = 50843.25
$ws.Range("M84").Value = -8001
$ws.Range("N84").Value = -62075.25
$ws.Range("H92").Value = 412
$ws.Range("I92").Value = 386.83334
$ws.Range("J92").Value = 449.75
$ws.Range("K92").Value = 1160.50002
$ws.Range("L92").Value = 1349.25
$ws.Range("M92").Value = 87.49998000000005
$ws.Range("N92").Value = -3845.25
$ws.Range("H113").Value = 588.94446
$ws.Range("J113").Value = 583.7273
$ws.Range("L113").Value = 1751.1819
$ws.Range("N113").Value = -6091.1819
$ws.Range("H132").Value = 1552
$ws.Range("I132").Value = 1412.4
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 12711.6
$ws.Range("L132").Value = 20250
$ws.Range("M132").Value = -10181.6
$ws.Range("N132").Value = -25310
$ws.Range("H137").Value = 3654
$ws.Range("I137").Value = 3028
$ws.Range("J137").Value = 4405.2
$ws.Range("K137").Value = 9084
$ws.Range("L137").Value = 13215.6
$ws.Range("M137").Value = -3984
$ws.Range("N137").Value = -23415.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 8996
$ws.Range("I21").Value = 6989
$ws.Range("K21").Value = 6989
$ws.Range("M21").Value = -6816
$ws.Range("H30").Value = 8996
$ws.Range("I30").Value = 6989
$ws.Range("K30").Value = 6989
$ws.Range("M30").Value = -6884
$ws.Range("H107").Value = 729.2857
$ws.Range("I107").Value = 461.7
$ws.Range("J107").Value = 1398.25
$ws.Range("K107").Value = 461.7
$ws.Range("L107").Value = 1398.25
$ws.Range("M107").Value = 1458.3
$ws.Range("N107").Value = -5238.25
$ws.Range("H122").Value = 3144.9062
$ws.Range("I122").Value = 2887.2083
$ws.Range("K122").Value = 8661.624899999999
$ws.Range("M122").Value = -6211.624899999999
$ws.Range("H132").Value = 4414.5293
$ws.Range("I132").Value = 3375.7693
$ws.Range("K132").Value = 10127.3079
$ws.Range("M132").Value = -7597.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 4354.2856
$ws.Range("I13").Value = 4974.7827
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 4974.7827
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -4834.7827
$ws.Range("N13").Value = -1780
$ws.Range("H61").Value = 3992.2
$ws.Range("I61").Value = 3855.3157
$ws.Range("J61").Value = 4425.6665
$ws.Range("K61").Value = 3855.3157
$ws.Range("L61").Value = 4425.6665
$ws.Range("M61").Value = -3653.3157
$ws.Range("N61").Value = -4829.6665
$ws.Range("H113").Value = 3992.2
$ws.Range("I113").Value = 3855.3157
$ws.Range("J113").Value = 4425.6665
$ws.Range("K113").Value = 3855.3157
$ws.Range("L113").Value = 4425.6665
$ws.Range("M113").Value = -1685.3157
$ws.Range("N113").Value = -8765.666499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 23733.334
$ws.Range("I58").Value = 19500
$ws.Range("K58").Value = 19500
$ws.Range("M58").Value = -19192
$ws.Range("H122").Value = 3017.375
$ws.Range("I122").Value = 3030
$ws.Range("J122").Value = 2996.3333
$ws.Range("K122").Value = 9090
$ws.Range("L122").Value = 8988.999899999999
$ws.Range("M122").Value = -6640
$ws.Range("N122").Value = -13888.9999
$ws.Range("H136").Value = 2967.7334
$ws.Range("J136").Value = 3200.4
$ws.Range("L136").Value = 9601.200000000001
$ws.Range("N136").Value = -14701.2
